# Adiciona "Caixa 43" (Organização / Cabides / Suíte / Baixa) como nova linha 350
# na planilha "Catálogo de Mudança", empurrando as linhas seguintes para baixo,
# e atualiza os totais na planilha "Resumo".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Catálogo de Mudança")

# Insere uma nova linha na posição 350 (as linhas 350+ descem para 351+)
$ws.Rows.Item(350).Insert()

# Copia a formatação "padrão" (estilo das colunas A-D e F) da linha acima,
# que já ocupa a linha 349 e permanece inalterada.
$ws.Range("A349:D349").Copy() | Out-Null
$ws.Range("A350:D350").PasteSpecial(-4122) | Out-Null

$ws.Range("F349").Copy() | Out-Null
$ws.Range("F350").PasteSpecial(-4122) | Out-Null

# Copia a formatação de prioridade "Baixa" (preenchimento cinza) da célula E12,
# que já usa esse estilo, para a nova célula E350.
$ws.Range("E12").Copy() | Out-Null
$ws.Range("E350").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Preenche os valores da nova linha (Caixa 43)
$ws.Cells.Item(350, 1).Value = "Caixa 43"
$ws.Cells.Item(350, 2).Value = "Organização"
$ws.Cells.Item(350, 3).Value = "Cabides"
$ws.Cells.Item(350, 4).Value = "Suíte"
$ws.Cells.Item(350, 5).Value = "Baixa"

# Atualiza o resumo (planilha "Resumo"): total de itens e contagem por prioridade "Baixa"
$resumo = $wb.Worksheets.Item("Resumo")
$resumo.Range("B3").Value = 455
$resumo.Range("B8").Value = 70
